$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Password" column from C to D to make room for the new
# "Nama Lengkap" (Full Name) column, then fill in the new layout.
$ws.Range("D1").Value = "Password"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D2").Value = 123456
$ws.Range("D3").Value = 123456

# New "Nama Lengkap" header + values in column C.
$ws.Range("C1").Value = "Nama Lengkap"
$ws.Range("C2").Value = "User Satu"
$ws.Range("C3").Value = "User Dua"

# "Jenis Pengguna" values switch from the long code to short role codes.
$ws.Range("A2").Value = "ADM"
$ws.Range("A3").Value = "DSN"

# Match the new Password column's width to column C's width.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Restore the active selection used after the edit.
$ws.Range("B5").Select() | Out-Null
